{"js": "async (context) => {\n  const body = context.document.body;\n\n  // --- 1. Rewrite the \"ways you will utilize the information\" paragraph ---\n  // The old text described a generic library-impact study; the new text\n  // describes the AUMAR augmented-reality study instead.\n  const oldUsageText =\n    \"The data that will be gathered within this survey is anonymous and does \" +\n    \"not contain any personal information. It will be used for the \" +\n    \"educational study targeting to reveal a modern library impact on the \" +\n    \"educational process via facts and connections (using data factsheet \" +\n    \"and/or database reports). As result, the obtained research results are \" +\n    \"planned to be published in a high-ranking journal.\";\n\n  const newUsageText =\n    \"The data that will be gathered within this survey is anonymous and does \" +\n    \"not contain any personal information. The aim of the proposed survey is \" +\n    \"to analyse the impact of augmented reality application in facilitating \" +\n    \"student orientation of the campus as well as the labs. The obtained \" +\n    \"research results are planned to be published in a high-ranking journal.\";\n\n  const usageResults = body.search(oldUsageText, { matchCase: true });\n  usageResults.load(\"text\");\n  await context.sync();\n\n  if (usageResults.items.length > 0) {\n    usageResults.items[0].insertText(newUsageText, \"Replace\");\n    await context.sync();\n  }\n\n  // --- 2. Tighten \"Comments :\" -> \"Comments:\" (only the first occurrence,\n  // in the Research Office assessment section; the Dean/President sections\n  // keep their original spacing). ---\n  const commentsResults = body.search(\"Comments :\", { matchCase: true });\n  commentsResults.load(\"text\");\n  await context.sync();\n\n  if (commentsResults.items.length > 0) {\n    commentsResults.items[0].insertText(\"Comments:\", \"Replace\");\n    await context.sync();\n  }\n};\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Rewrite the \"ways you will utilize the information\" paragraph ---\n# The old text described a generic library-impact study; the new text\n# describes the AUMAR augmented-reality study instead.\n$oldUsageText = \"The data that will be gathered within this survey is anonymous and does not contain any personal information. It will be used for the educational study targeting to reveal a modern library impact on the educational process via facts and connections (using data factsheet and/or database reports). As result, the obtained research results are planned to be published in a high-ranking journal.\"\n$newUsageText = \"The data that will be gathered within this survey is anonymous and does not contain any personal information. The aim of the proposed survey is to analyse the impact of augmented reality application in facilitating student orientation of the campus as well as the labs. The obtained research results are planned to be published in a high-ranking journal.\"\n\n$usageRange = $d.Content\n$usageRange.Find.ClearFormatting()\n$usageRange.Find.Replacement.ClearFormatting()\n$usageRange.Find.Execute($oldUsageText, $false, $false, $false, $false, $false, $true, 1, $false, $newUsageText, 1)\n\n# --- 2. Tighten \"Comments :\" -> \"Comments:\" (only the first occurrence, in\n# the Research Office assessment section; the Dean/President sections keep\n# their original spacing). ---\n$commentsRange = $d.Content\n$commentsRange.Find.ClearFormatting()\n$commentsRange.Find.Replacement.ClearFormatting()\n$commentsRange.Find.Execute(\"Comments :\", $false, $false, $false, $false, $false, $true, 1, $false, \"Comments:\", 1)\n"}
